$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 5) mirroring the formatting of rows 3/4
$ws.Range("A5").Value = "W2"

$values5 = @(9.89, 9.09, 10.35, 12.18, 10.5, 9.09, 9.9499999999999993, 11.3, 10.15, 10.11, 12.35, 11.11)
for ($i = 0; $i -lt $values5.Length; $i++) {
    $col = $i + 2  # B=2
    $ws.Cells.Item(5, $col).Value = $values5[$i]
}

# Copy style/border formatting from row 4 (A4:M4) to row 5 (A5:M5)
$ws.Range("A4:M4").Copy()
$ws.Range("A5:M5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to match the target state
$ws.Range("F9").Select()
